$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-16 Tuesday", "2025-09-17 Wednesday"),
    @("978÷5=", "778÷7="),
    @("289÷9=", "768÷9="),
    @("518÷7=", "216÷4="),
    @("352÷8=", "456÷5="),
    @("531÷7=", "733÷8="),
    @("551÷7=", "976÷4="),
    @("588÷5=", "926÷2="),
    @("117÷7=", "425÷9="),
    @("590÷6=", "952÷5="),
    @("687÷4=", "909÷7="),
    @("894÷2=", "700÷9="),
    @("953÷4=", "775÷6="),
    @("154÷3=", "990÷5="),
    @("826÷6=", "938÷2="),
    @("187÷2=", "281÷5="),
    @("982÷8=", "522÷9="),
    @("275÷7=", "570÷7="),
    @("696÷9=", "805÷6="),
    @("186÷6=", "860÷4="),
    @("740÷2=", "888÷9="),
    @("638÷4=", "572÷9="),
    @("276÷8=", "127÷9="),
    @("829÷3=", "703÷9="),
    @("623÷2=", "719÷4="),
    @("341÷7=", "157÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
